# Updates cryptos list values (Price / Volume(1h) columns, plus a few
# Coin/Link cells in rows 38-41 where row order shifted) to match the
# latest scrape, per commit "Updated cryptos list on Sat Oct 12 18:59:40 UTC 2024 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers formatted as plain text (e.g. thousands
# separated with '.', or using subscript-zero notation); setting .Value directly
# would let Excel silently reinterpret them as real numbers and reformat/round
# them, so we force text format first, write the value, then drop the explicit
# format again so the cell keeps using the original (unstyled) look.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '62.984.54'
$ws.Range('E2').Value = '  +0.15%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.470.91'
$ws.Range('E3').Value = '  +1.05%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
Set-TextValue $ws.Range('D5') '576.40'
$ws.Range('E5').Value = '  -0.76%  '

# Row 6
Set-TextValue $ws.Range('D6') '146.52'
$ws.Range('E6').Value = '  +0.46%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('E8').Value = '  -0.13%  '

# Row 9
Set-TextValue $ws.Range('D9') '2.470.53'
$ws.Range('E9').Value = '  +1.10%  '

# Row 10
$ws.Range('E10').Value = '  +0.56%  '

# Row 11
$ws.Range('E11').Value = '  +1.20%  '

# Row 12
$ws.Range('E12').Value = '  +0.61%  '

# Row 13
$ws.Range('E13').Value = '  +0.11%  '

# Row 14
Set-TextValue $ws.Range('D14') '28.94'
$ws.Range('E14').Value = '  +7.14%  '

# Row 15
$ws.Range('E15').Value = '  -0.20%  '

# Row 16
$ws.Range('E16').Value = '  +1.00%  '

# Row 17
Set-TextValue $ws.Range('D17') '63.028.18'
$ws.Range('E17').Value = '  +0.55%  '

# Row 18
Set-TextValue $ws.Range('D18') '2.467.98'
$ws.Range('E18').Value = '  +1.24%  '

# Row 19
Set-TextValue $ws.Range('D19') '8.19'
$ws.Range('E19').Value = '  +2.99%  '

# Row 20
Set-TextValue $ws.Range('D20') '11.03'
$ws.Range('E20').Value = '  +0.92%  '

# Row 21
Set-TextValue $ws.Range('D21') '329.37'
$ws.Range('E21').Value = '  +0.88%  '

# Row 23
$ws.Range('E23').Value = '  +7.64%  '

# Row 24
$ws.Range('E24').Value = '  -0.10%  '

# Row 25
Set-TextValue $ws.Range('D25') '66.28'
$ws.Range('E25').Value = '  +0.79%  '

# Row 26
Set-TextValue $ws.Range('D26') '666.91'
$ws.Range('E26').Value = '  +8.27%  '

# Row 27
Set-TextValue $ws.Range('D27') '9.52'
$ws.Range('E27').Value = '  +13.73%  '

# Row 28
Set-TextValue $ws.Range('D28') '0.0₃0985'
$ws.Range('E28').Value = '  +0.48%  '

# Row 30
Set-TextValue $ws.Range('D30') '0.998'
$ws.Range('E30').Value = '  +957.87%  '

# Row 31
$ws.Range('E31').Value = '  +2.42%  '

# Row 32
$ws.Range('E32').Value = '  -0.98%  '

# Row 33
$ws.Range('E33').Value = '  -0.90%  '

# Row 34
$ws.Range('E34').Value = '  -2.68%  '

# Row 35
$ws.Range('E35').Value = '  +3.71%  '

# Row 36
$ws.Range('E36').Value = '  +0.00%  '

# Row 37
$ws.Range('E37').Value = '  +0.69%  '

# Row 38
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D38') '152.62'
$ws.Range('E38').Value = '  +0.01%  '

# Row 39
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range('D39') '0.372'
$ws.Range('E39').Value = '  -0.22%  '

# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D40') '5.42'
$ws.Range('E40').Value = '  +0.77%  '

# Row 41
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D41') '18.73'
$ws.Range('E41').Value = '  +0.72%  '

# Row 42
Set-TextValue $ws.Range('D42') '2.72'
$ws.Range('E42').Value = '  -0.63%  '

# Row 43
$ws.Range('E43').Value = '  -0.64%  '

# Row 45
$ws.Range('E45').Value = '  -5.97%  '

# Row 46
Set-TextValue $ws.Range('D46') '150.77'
$ws.Range('E46').Value = '  +4.72%  '

# Row 47
Set-TextValue $ws.Range('D47') '15.14'
$ws.Range('E47').Value = '  +26.74%  '

# Row 48
$ws.Range('E48').Value = '  +0.75%  '

# Row 49
Set-TextValue $ws.Range('D49') '20.67'
$ws.Range('E49').Value = '  +2.57%  '

# Row 50
$ws.Range('E50').Value = '  +1.36%  '

# Row 51
$ws.Range('E51').Value = '  -0.14%  '
